# Updates cryptos list (Price / Volume(1h) columns) with refreshed values,
# and reorders a few coin rows (FraxShare/Aave, EnergySwap/SynthetixNetwork/Aptos)
# to match the latest scrape, as published on Wed Jul 26 21:09:08 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking price strings (e.g. "29.610.98", "0.9995")
# that must stay as literal text (matching the original inline-string cells) instead
# of being auto-converted to numbers by Excel. Forcing the Text number format first
# ensures the assigned Value is stored as a string.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.610.98"
$ws.Range("E2").Value = "  +1.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.881.39"
$ws.Range("E3").Value = "  +1.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7271"
$ws.Range("E5").Value = "  +3.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.83"
$ws.Range("E6").Value = "  +0.80%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9997"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07894"
$ws.Range("E8").Value = "  -2.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3089"
$ws.Range("E9").Value = "  +1.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.27"
$ws.Range("E10").Value = "  +8.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08245"
$ws.Range("E11").Value = "  +1.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.878.57"
$ws.Range("E12").Value = "  +1.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7291"
$ws.Range("E13").Value = "  +2.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.283"
$ws.Range("E14").Value = "  +2.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.76"
$ws.Range("E15").Value = "  +1.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.623.11"
$ws.Range("E16").Value = "  +1.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.876"
$ws.Range("E17").Value = "  +1.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "243.94"
$ws.Range("E18").Value = "  +3.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007892"
$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.44"
$ws.Range("E20").Value = "  +0.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.124.13"
$ws.Range("E21").Value = "  +0.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9989"
$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.791"
$ws.Range("E24").Value = "  +5.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1559"
$ws.Range("E25").Value = "  +8.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.23"
$ws.Range("E26").Value = "  +0.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.007"
$ws.Range("E27").Value = "  +0.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.44"
$ws.Range("E28").Value = "  +1.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.966"
$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.370"
$ws.Range("E30").Value = "  -4.27%  "

$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.367"
$ws.Range("E32").Value = "  -0.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.113"
$ws.Range("E33").Value = "  +1.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05268"
$ws.Range("E34").Value = "  +1.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.203"
$ws.Range("E35").Value = "  +2.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7190"
$ws.Range("E36").Value = "  +1.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  +0.28%  "

$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01869"
$ws.Range("E39").Value = "  +1.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.718"
$ws.Range("E40").Value = "  -0.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.181.09"
$ws.Range("E41").Value = "  +3.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9125"
$ws.Range("E42").Value = "  -0.72%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "6.012"
$ws.Range("E43").Value = "  +2.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "72.06"
$ws.Range("E44").Value = "  +2.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4326"
$ws.Range("E45").Value = "  +1.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.82"
$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5369"
$ws.Range("E48").Value = "  -0.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.785"
$ws.Range("E49").Value = "  +0.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").Value = "2.887"
$ws.Range("E50").Value = "  +5.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "9.248"
$ws.Range("E51").Value = "  +0.56%  "
